$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Phase 1: row 2 carries a sticky row-level "customFormat" flag from
# the original sheet that plain cell-level operations cannot remove.
# Insert a fresh row in its place and delete the old (tainted) row
# that gets pushed down - this nets out to zero row-count change.
# ------------------------------------------------------------------
$ws.Rows(2).Insert()
$ws.Rows(3).Delete()

# ------------------------------------------------------------------
# Phase 2: stash copies of every cell style we will need (far away,
# so later clears/merges on A1:B14 do not disturb them).
# ------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("D30").PasteSpecial(-4122)            # style 6: header (blue/navy, align left)

$ws.Range("A3").Copy()
$ws.Range("D31").PasteSpecial(-4122)            # style 3: red bold, fill2, centered

$ws.Range("B3").Copy()
$ws.Range("D32").PasteSpecial(-4122)            # style 4: red bold, fill2

$ws.Range("B4").Copy()
$ws.Range("D33").PasteSpecial(-4122)            # style 5: green bold, fill2

$ws.Range("D31").Copy()
$ws.Range("D34").PasteSpecial(-4122)
$ws.Range("D34").Font.Color = 65280             # new style 7: green bold, fill2, centered

$ws.Range("D31").Copy()
$ws.Range("D35").PasteSpecial(-4122)
$ws.Range("D35").Font.Color = 65535             # new style 8: yellow bold, fill2, centered

$ws.Range("D32").Copy()
$ws.Range("D36").PasteSpecial(-4122)
$ws.Range("D36").Font.Color = 65535             # new style 9: yellow bold, fill2

# ------------------------------------------------------------------
# Phase 3: wipe the old table content (values, styles).  A1:B1 is
# already merged in the source and stays merged in the target, so it
# is intentionally left alone (re-merging needlessly allocates a
# spurious extra style entry).
# ------------------------------------------------------------------
$ws.Range("A1:B14").Clear()

# ------------------------------------------------------------------
# Phase 4: rebuild the table with the new content.
# ------------------------------------------------------------------
$ws.Range("D30").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = "Autonomous Modes"
$ws.Range("D30").Copy()
$ws.Range("B1").PasteSpecial(-4122)

$ws.Range("D31").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = "MODE #"
$ws.Range("D32").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("B2").Value = "Defense"

$ws.Range("D34").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 8
$ws.Range("D33").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B3").Value = "Corner Shot"

$ws.Range("D34").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 100
$ws.Range("D33").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B4").Value = "Low Bar One Ball (w Gyro)"

$ws.Range("D34").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 101
$ws.Range("D33").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B5").Value = "Portcullis One Ball (w Gyro)"

$ws.Range("D34").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = 102
$ws.Range("D33").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("B6").Value = "Cheval One Ball (w Gyro)"

$ws.Range("D34").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = 103
$ws.Range("D33").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B7").Value = "Rough Terrain One Ball (w Gyro)"

$ws.Range("D35").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = 200
$ws.Range("D36").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B8").Value = "Low Bar Two Ball  w Spybot (w Gyro)"

$ws.Range("D34").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = "default"
$ws.Range("D33").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").Value = "Corner Shot"

# ------------------------------------------------------------------
# Phase 5: rows 10-14 no longer hold any content - drop them so the
# sheet dimension shrinks back down to A1:B9.
# ------------------------------------------------------------------
$ws.Rows("10:14").Delete()

# ------------------------------------------------------------------
# Phase 6: drop the scratch cells.
# ------------------------------------------------------------------
$ws.Range("D30:D36").Delete()

$ws.Range("A1").Select()
